$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E4: ContentType.JSON -> None
$ws.Range("E4").Value = "None"

# Rows 7 and 8 effectively swap their API-test data (endpoint, method,
# path params, auth type/creds, expected status, assertions), with the
# auth credentials text on row 8 updated to "username:password".

# Row 7 (was /register POST ... BasicAuth/username and password/400/None)
# becomes the GET /users/{id} test case.
$ws.Range("C7").Value = "/users/{id}"
$ws.Range("D7").Value = "GET"
$ws.Range("H7").Value = "id=3"
$ws.Range("J7").Value = "None"
$ws.Range("K7").Value = "None"
$ws.Range("L7").Value = "'200"
$ws.Range("N7").Value = "data.last_name: 'Wong'"

# Row 8 (was /users/{id} GET ... id=3/None/None/200/data.last_name) becomes
# the POST /register test case.
$ws.Range("C8").Value = "/register"
$ws.Range("D8").Value = "POST"
$ws.Range("H8").Value = "None"
$ws.Range("J8").Value = "BasicAuth"
$ws.Range("K8").Value = "username:password"
$ws.Range("L8").Value = "'400"
$ws.Range("N8").Value = "None"
